$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = 5.210200000000001
$ws.Range("B12").Value = 5.162999999999999
$ws.Range("B18").Value = 6.990499999999997
$ws.Range("B37").Value = 8.746400000000001
$ws.Range("B55").Value = 6.291099999999995
$ws.Range("B68").Value = 4.871399999999997
$ws.Range("B77").Value = 9.229200000000008
$ws.Range("B78").Value = 9.582000000000003
